# Actualización automática del tracker

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Results came in for a few pending matches (rows 7, 8, 26)
$ws.Range("G7").Value = "Acierto"
$ws.Range("H7").Value = 0.91

$ws.Range("G8").Value = "Fallo"
$ws.Range("H8").Value = -1

$ws.Range("G26").Value = "Fallo"
$ws.Range("H26").Value = -1

# Normalize event_id column (A9:A27) to real numbers instead of text
$eventIdRows = @(9, 10, 11, 12, 13, 14, 15, 16, 17, 18, 19, 20, 21, 22, 23, 24, 25, 26, 27)
foreach ($r in $eventIdRows) {
    $cell = $ws.Cells.Item($r, 1)
    $cell.Value = [double]$cell.Value2
}
